# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "103.51") need the
# columns underlying number format forced to Text first, otherwise Excel
# auto-converts the assigned string into a floating point number (losing
# the original formatted precision, e.g. "8.20" -> 8.2).
$forceTextCells = @("D5", "D6", "D7", "D10", "D11", "D12", "D14", "D15", "D19", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D30", "D31", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D42", "D43", "D44", "D45", "D48", "D49")
foreach ($cellRef in $forceTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '43.771.96'
$ws.Range("E2").Value = '  -0.09%  '
$ws.Range("D3").Value = '2.291.80'
$ws.Range("E3").Value = '  -1.12%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '103.51'
$ws.Range("E5").Value = '  +6.45%  '
$ws.Range("D6").Value = '270.67'
$ws.Range("E6").Value = '  -0.73%  '
$ws.Range("D7").Value = '0.625'
$ws.Range("E7").Value = '  -0.53%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E9").Value = '  -2.77%  '
$ws.Range("D10").Value = '45.98'
$ws.Range("E10").Value = '  +1.45%  '
$ws.Range("D11").Value = '0.0935'
$ws.Range("E11").Value = '  -1.75%  '
$ws.Range("D12").Value = '8.20'
$ws.Range("E12").Value = '  +2.16%  '
$ws.Range("E13").Value = '  +1.69%  '
$ws.Range("D14").Value = '15.53'
$ws.Range("E14").Value = '  +0.04%  '
$ws.Range("D15").Value = '0.856'
$ws.Range("E15").Value = '  -2.37%  '
$ws.Range("D16").Value = '2.291.52'
$ws.Range("E16").Value = '  -1.26%  '
$ws.Range("D17").Value = '43.789.52'
$ws.Range("E17").Value = '  +0.05%  '
$ws.Range("E18").Value = '  +0.35%  '
$ws.Range("D19").Value = '6.27'
$ws.Range("D20").Value = '72.24'
$ws.Range("E20").Value = '  -1.55%  '
$ws.Range("D21").Value = '2.51'
$ws.Range("E21").Value = '  +10.44%  '
$ws.Range("D22").Value = '233.58'
$ws.Range("E22").Value = '  -2.62%  '
$ws.Range("D23").Value = '2.91'
$ws.Range("E23").Value = '  +14.25%  '
$ws.Range("D24").Value = '9.30'
$ws.Range("E24").Value = '  -1.16%  '
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("D26").Value = '11.31'
$ws.Range("E26").Value = '  -0.57%  '
$ws.Range("D27").Value = '40.49'
$ws.Range("E27").Value = '  +5.91%  '
$ws.Range("E28").Value = '  -1.66%  '
$ws.Range("E29").Value = '  -3.23%  '
$ws.Range("D30").Value = '177.78'
$ws.Range("E30").Value = '  +1.50%  '
$ws.Range("D31").Value = '21.82'
$ws.Range("E31").Value = '  -2.70%  '
$ws.Range("E32").Value = '  -1.58%  '
$ws.Range("D33").Value = '5.54'
$ws.Range("E33").Value = '  +0.76%  '
$ws.Range("D34").Value = '4.89'
$ws.Range("E34").Value = '  +10.03%  '
$ws.Range("D35").Value = '0.127'
$ws.Range("E35").Value = '  -0.48%  '
$ws.Range("D36").Value = '0.111'
$ws.Range("E36").Value = '  +1.49%  '
$ws.Range("E37").Value = '  -0.97%  '
$ws.Range("D38").Value = '3.56'
$ws.Range("E38").Value = '  +4.88%  '
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").Value = '2.33'
$ws.Range("E39").Value = '  -3.25%  '
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").Value = '0.236'
$ws.Range("E40").Value = '  -4.01%  '
$ws.Range("E41").Value = '  -2.70%  '
$ws.Range("D42").Value = '65.50'
$ws.Range("E42").Value = '  +4.37%  '
$ws.Range("D43").Value = '12.19'
$ws.Range("E43").Value = '  -0.99%  '
$ws.Range("D44").Value = '5.31'
$ws.Range("E44").Value = '  -0.86%  '
$ws.Range("D45").Value = '8.80'
$ws.Range("E45").Value = '  -4.71%  '
$ws.Range("E46").Value = '  -1.38%  '
$ws.Range("E47").Value = '  +2.27%  '
$ws.Range("D48").Value = '99.29'
$ws.Range("E48").Value = '  -1.19%  '
$ws.Range("D49").Value = '1.55'
$ws.Range("E49").Value = '  +11.65%  '
$ws.Range("E50").Value = '  +4.36%  '
$ws.Range("D51").Value = '2.525.23'
$ws.Range("E51").Value = '  -0.76%  '

# Restore the default "Normal" style on the cells we force-formatted as text
# so no stray number-format style is left attached to them.
foreach ($cellRef in $forceTextCells) {
    $ws.Range($cellRef).Style = "Normal"
}
